$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header cells for the new columns, matching the bold/border style used by other headers (style index 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # AD
    $ws.Cells.Item($r, 31).Value = 89   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
